$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Clear out the whole previously-used area so we can rebuild it
# cleanly (content + formatting) to match the new layout.
# ------------------------------------------------------------------
$ws.Range("A1:E18").Clear()

# ------------------------------------------------------------------
# Header row
# ------------------------------------------------------------------
$ws.Range("B1").Value2 = "Start"
$ws.Range("C1").Value2 = "End"
$ws.Range("D1").Value2 = "Duration"
$ws.Range("E1").Value2 = "Notes"

# ------------------------------------------------------------------
# Time-log data rows 2-14 (B=Start, C=End, D=Duration formula, E=Notes)
# ------------------------------------------------------------------
$ws.Range("B2").Value2 = 0.07986111111111110494
$ws.Range("B2").NumberFormat = "h:mm"
$ws.Range("C2").Value2 = 0.13541666666666665741
$ws.Range("C2").NumberFormat = "h:mm"
$ws.Range("D2").Formula = "=C2-B2"
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("E2").Value2 = "Updating notes in new quarto document. Reading up on alternative clustering procedures and methods."

$ws.Range("B3").Value2 = 0.10763888888888889506
$ws.Range("B3").NumberFormat = "h:mm"
$ws.Range("C3").Value2 = 0.13888888888888889506
$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("D3").Formula = "=C3-B3"
$ws.Range("D3").NumberFormat = "h:mm"

$ws.Range("B4").Value2 = 0.30208333333333331483
$ws.Range("B4").NumberFormat = "h:mm"
$ws.Range("C4").Value2 = 0.31944444444444447528
$ws.Range("C4").NumberFormat = "h:mm"
$ws.Range("D4").Formula = "=C4-B4"
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("E4").Value2 = "Reading more on cluster methods and beginning work on website. "

$ws.Range("B5").Value2 = 0.42361111111111110494
$ws.Range("B5").NumberFormat = "h:mm"
$ws.Range("C5").Value2 = 0.50694444444444441977
$ws.Range("C5").NumberFormat = "h:mm"
$ws.Range("D5").Formula = "=C5-B5"
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("E5").Value2 = "Working up website to show updates to clustering procedures. "

$ws.Range("B6").Value2 = 0.44444444444444441977
$ws.Range("B6").NumberFormat = "h:mm"
$ws.Range("C6").Value2 = 0.5090277777777777457
$ws.Range("C6").NumberFormat = "h:mm"
$ws.Range("D6").Formula = "=C6-B6"
$ws.Range("D6").NumberFormat = "h:mm"
$ws.Range("E6").Value2 = "Research factor analysis and principal components analysis methods. Exploratory coding work to look at viability of these methods."

$ws.Range("B7").Value2 = 0.10416666666666667129
$ws.Range("B7").NumberFormat = "h:mm"
$ws.Range("C7").Value2 = 0.13888888888888889506
$ws.Range("C7").NumberFormat = "h:mm"
$ws.Range("D7").Formula = "=C7-B7"
$ws.Range("D7").NumberFormat = "h:mm"
$ws.Range("E7").Value2 = "Selecting new variables for principal components analysis and revisting data cleaning portion"

$ws.Range("B8").Value2 = 0.31944444444444447528
$ws.Range("B8").NumberFormat = "h:mm"
$ws.Range("C8").Value2 = 0.33333333333333331483
$ws.Range("C8").NumberFormat = "h:mm"
$ws.Range("D8").Formula = "=C8-B8"
$ws.Range("D8").NumberFormat = "h:mm"
$ws.Range("E8").Value2 = "Coding principal components scores."

$ws.Range("B9").Value2 = 0.35694444444444445308
$ws.Range("B9").NumberFormat = "h:mm"
$ws.Range("C9").Value2 = 0.4027777777777777346
$ws.Range("C9").NumberFormat = "h:mm"
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("D9").NumberFormat = "h:mm"
$ws.Range("E9").Value2 = "Generate clusters using PC scores and quantile groupings. Visual inspection of different combinations of variables to see if we can make intuitive sense of the different groupings and combinations."

$ws.Range("B10").Value2 = 0.45833333333333331483
$ws.Range("B10").NumberFormat = "h:mm"
$ws.Range("C10").Value2 = 0.4722222222222222654
$ws.Range("C10").NumberFormat = "h:mm"
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("D10").NumberFormat = "h:mm"
$ws.Range("E10").Value2 = "Working on a table that will show the the principal component scores across the three key dimensiosn broken down by percentiles. I think we need to break out the demographic into a race and economic dimension, though, as major metro areas like NY are occupying a stange middle ground in the current iteration. I think some things some signs of ""weak"" economic performance, like unemployment, are working against other more positive indicators, like total population and diversity. "

$ws.Range("B11").Value2 = 0.16666666666666665741
$ws.Range("B11").NumberFormat = "h:mm"
$ws.Range("C11").Value2 = 0.2048611111111111327
$ws.Range("C11").NumberFormat = "h:mm"
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("D11").NumberFormat = "h:mm"
$ws.Range("E11").Value2 = "Broke the demographic variable out into two dimensions. Tried different cluster techniques on the raw PC scores and percentile PC scores. Tried various visualization techniques to try to show the different dimensions captured by the PC scores."

$ws.Range("B12").Value2 = 0.4027777777777777346
$ws.Range("B12").NumberFormat = "h:mm"
$ws.Range("C12").Value2 = 0.4375
$ws.Range("C12").NumberFormat = "h:mm"
$ws.Range("D12").Formula = "=C12-B12"
$ws.Range("D12").NumberFormat = "h:mm"

$ws.Range("B13").Value2 = 0.16666666666666665741
$ws.Range("B13").NumberFormat = "h:mm"
$ws.Range("C13").Value2 = 0.20833333333333334259
$ws.Range("C13").NumberFormat = "h:mm"
$ws.Range("D13").Formula = "=C13-B13"
$ws.Range("D13").NumberFormat = "h:mm"

$ws.Range("B14").Value2 = 0.125
$ws.Range("B14").NumberFormat = "h:mm"
$ws.Range("C14").Value2 = 0.13888888888888889506
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("D14").Formula = "=C14-B14"
$ws.Range("D14").NumberFormat = "h:mm"

# ------------------------------------------------------------------
# Blank spacer rows 15-19 (kept as time-formatted blank cells, like
# the original template rows below the log)
# ------------------------------------------------------------------
foreach ($r in 15,16,17) {
  $ws.Range("B$r").NumberFormat = "h:mm"
  $ws.Range("C$r").NumberFormat = "h:mm"
  $ws.Range("D$r").NumberFormat = "h:mm"
}
$ws.Range("D18").NumberFormat = "h:mm"
$ws.Range("D19").NumberFormat = "h:mm"

# ------------------------------------------------------------------
# Totals block (rows 20-22)
# ------------------------------------------------------------------

# Values/formulas first
$ws.Range("A20").Value2 = "Total Time"
$ws.Range("D20").Formula = "=SUM(D2:D14)"
$ws.Range("A21").Value2 = "Rate"
$ws.Range("D21").Value2 = 200
$ws.Range("A22").Value2 = "Total Billables"
$ws.Range("D22").Formula = "=D20*24*D21"

# Number formats, in the same order the style table expects them
$ws.Range("D22").NumberFormat = "#,##0.00;[Red]#,##0.00"
$ws.Range("D21").NumberFormat = """$""#,##0;[Red]""$""#,##0"

# Double top border across the Total Time row (A:E)
$ws.Range("B20:E20").Borders.Item(8).LineStyle = -4119
$ws.Range("D20").NumberFormat = "h:mm"
$ws.Range("A20").Borders.Item(8).LineStyle = -4119
$ws.Range("A20").Font.Bold = $true

# Bold labels
$ws.Range("A21").Font.Bold = $true
$ws.Range("A22").Font.Bold = $true

# ------------------------------------------------------------------
# Selection / view state
# ------------------------------------------------------------------
$ws.Range("C15").Select()
